# Refined metadata to be additional tab
#
# 1. Refresh the "time_taken" query timestamps on the existing "data" sheet.
# 2. Add a new "metadata" worksheet (after "data") describing the panel
#    query itself (name/id/version/retrieval time/request url).

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- 1. Updated timestamps on the "data" sheet (column F, rows 2-9) ---
$data.Range("F2").Value = "2021-10-05 14:20:52.531750"
$data.Range("F3").Value = "2021-10-05 14:20:52.531757"
$data.Range("F4").Value = "2021-10-05 14:20:52.531760"
$data.Range("F5").Value = "2021-10-05 14:20:52.531763"
$data.Range("F6").Value = "2021-10-05 14:20:52.531766"
$data.Range("F7").Value = "2021-10-05 14:20:52.531768"
$data.Range("F8").Value = "2021-10-05 14:20:52.531770"
$data.Range("F9").Value = "2021-10-05 14:20:52.531773"

# --- 2. New "metadata" worksheet, inserted after "data" ---
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "metadata"

# Match the outline properties already used on the "data" sheet.
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

# Header row
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Reuse the exact same (bold / bordered / centered) cell style already used
# for the "data" sheet's header row and its "A" index column, by copying
# formats across rather than re-creating a new font/style from scratch.
$data.Range("B1:F1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)
$data.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data row
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Hyperthyroidism"
$ws.Range("C2").Value = 236

# Force "2.8" to be stored as literal text (not the number 2.8), matching
# the source data, then drop the temporary text number-format so the cell
# carries no explicit style (same as the rest of the un-styled data row).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2.8"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = "2021-08-24T12:23:21.327542Z"
$ws.Range("F2").Value = "2021-10-05 14:20:52.528078"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/236/?format=json"

# Match page margins to the "data" sheet.
$ws.PageSetup.LeftMargin = $data.PageSetup.LeftMargin
$ws.PageSetup.RightMargin = $data.PageSetup.RightMargin
$ws.PageSetup.TopMargin = $data.PageSetup.TopMargin
$ws.PageSetup.BottomMargin = $data.PageSetup.BottomMargin
$ws.PageSetup.HeaderMargin = $data.PageSetup.HeaderMargin
$ws.PageSetup.FooterMargin = $data.PageSetup.FooterMargin

# Leave the original "data" sheet selected/active, as before the edit.
$data.Activate()
